$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (D, E, F) ---
$ws.Columns.Item(4).ColumnWidth = 15.17
$ws.Columns.Item(5).ColumnWidth = 14.17
$ws.Columns.Item(6).ColumnWidth = 15.17

# --- Unmerge F7:F21 before editing F6/F7, will re-merge as F6:F21 later ---
$ws.Range("F7:F21").UnMerge()

# --- Cell value updates ---
$ws.Range("A2").Value = "ו"
$ws.Range("C2").Value = "אבנר`nשגיא"
$ws.Range("D2").Value = "דימה`nשבצוב"
$ws.Range("E2").Value = "אור`nאסרף"
$ws.Range("F2").Value = "לוטם`nסיני"
$ws.Range("C5").Value = "אנזו`nשרעבי"
$ws.Range("D5").Value = "דימנטמן`nמטמוני"
$ws.Range("E5").Value = "דותן`nליאור"
$ws.Range("C8").Value = "אלכסיי`nלומיאנסקי"
$ws.Range("D8").Value = "דעאל`nלואיס"
$ws.Range("E8").Value = "אנדי`nדוד"
$ws.Range("C11").Value = "כלפה`nשמעון"
$ws.Range("D11").Value = "ארד`nיואל"
$ws.Range("E11").Value = "נפמן`nסדון"
$ws.Range("C14").Value = "מרדש`nמשה"
$ws.Range("D14").Value = "אסף`nדורון"
$ws.Range("E14").Value = "קריספין`nרווה"
$ws.Range("C17").Value = "דבוש`nפיאצה"
$ws.Range("D17").Value = "שראל`nשרעבי"
$ws.Range("E17").Value = "דימה`nשבצוב"
$ws.Range("C20").Value = "לוטם`nסיני"
$ws.Range("D20").Value = "אור`nאסרף"
$ws.Range("E20").Value = "אבנר`nאיתי כהן"
$ws.Range("F22").Value = "דעאל`nלומיאנסקי"
$ws.Range("C23").Value = "שגיא`nליאור"
$ws.Range("D23").Value = "אנדי`nדוד"
$ws.Range("E23").Value = "אנזו`nכלפה"
$ws.Range("A24").Value = "שבת"
$ws.Range("C26").Value = "ארד`nיואל"
$ws.Range("D26").Value = "מרדש`nמשה"
$ws.Range("E26").Value = "לואיס`nשמעון"
$ws.Range("F26").Value = "נפמן`nסדון"
$ws.Range("C29").Value = "אסף`nדורון"
$ws.Range("D29").Value = "קריספין`nרווה"
$ws.Range("E29").Value = "דבוש`nפיאצה"
$ws.Range("C32").Value = "שראל`nשרעבי"
$ws.Range("D32").Value = "דימה`nשבצוב"
$ws.Range("E32").Value = "אור`nאסרף"
$ws.Range("C35").Value = "אבנר`nלומיאנסקי"
$ws.Range("D35").Value = "איתי כהן`nכלפה"
$ws.Range("E35").Value = "לוטם`nסיני"
$ws.Range("C38").Value = "דעאל`nשגיא"
$ws.Range("D38").Value = "לואיס`nאנזו"
$ws.Range("E38").Value = "אנדי`nדוד"
$ws.Range("C41").Value = "נפמן`nסדון"
$ws.Range("D41").Value = "ארד`nיואל"
$ws.Range("E41").Value = "ליאור`nמשה"

# --- F6/F7 swap: F6 becomes a single space, F7 becomes empty ---
$ws.Range("F6").Value = " "
$ws.Range("F7").ClearContents()

# --- Re-merge F6:F21 ---
$ws.Range("F6:F21").Merge()

